$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time-range labels in column C (rows 2-7)
$ws.Range("C2").Value = "9:30-9:35"
$ws.Range("C3").Value = "9:35-9:40"
$ws.Range("C4").Value = "10:55-11:0"
$ws.Range("C5").Value = "11:0-11:5"
$ws.Range("C6").Value = "22:30-22:35"
$ws.Range("C7").Value = "22:35-22:40"

# Move the active selection to C12, matching the saved cursor position
$ws.Range("C12").Select()
